# ESTRATEGIA.docx edit script
# 1) "D"+"escripció"+"n de "+"las"+" " runs -> single run "Descripción de las "
# 2) "Se entendió que si..." paragraph rewritten to
#    "El paciente pertenecerá o no aun grupo familiar si al buscar..." with
#    proofErr spell-check markers around "aun", "ultimos", "digitos"
# 3) Paragraph with PROCEDURE SOLARIS.especialidadesMasBonosConsultaUsados:
#    move <w:lastRenderedPageBreak/> from before "los bonos..." to before
#    "PROCEDURE " and merge the two runs that used to be split by it
# 4) Paragraph with lone "Notas:" right after the REGISTRAR LLEGADA heading
#    gains a <w:lastRenderedPageBreak/> before its text

$d = $word.ActiveDocument

# --- Edit 1: merge "Descripcion de las " runs -------------------------------
$full = $d.Content
$full.Find.Execute("Descripción de las ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target1 = $d.Range($full.Start, $full.End)
$target1.Find.Execute("Descripción de las ", $false, $false, $false, $false, $false, $true, 1, $false, "Descripción de las ", 2) | Out-Null

# --- Edit 2: rewrite "Se entendio..." paragraph -----------------------------
$full = $d.Content
$full.Find.Execute("Se entendió que si el numero de familiares es 0 , entonces no pertenece a ningún grupo familiar", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target2 = $d.Range($full.Start, $full.End)

$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">El paciente pertenecerá o no </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>aun</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> grupo familiar si al buscar en la tabla pacientes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> encontramos más de 1 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">que tenga similar código de paciente y solo se diferencien en los 2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>ultimos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>digitos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target2.InsertXML($xml2)

# --- Edit 3: especialidadesMasBonosConsultaUsados paragraph ----------------
$full = $d.Content
$full.Find.Execute("PROCEDURE SOLARIS.especialidadesMasBonosConsultaUsados traerá las especialidades con mas bonos de consulta usados, aquí se tomo en cuenta solo los bonos que tiene asociado una consulta, se filtra que sean bonos de tipo consulta, utilizados y que la fecha de la consulta coincida con la fecha pedida", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target3 = $d.Range($full.Start, $full.End)

$xml3 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00451D68"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">PROCEDURE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00451D68"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>SOLARIS.especialidadesMasBonosConsultaUsados</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00451D68"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> traerá las especialidades con mas bonos de consulta usados, aquí se tomo en cuenta solo los bonos que tiene asociado una consulta, se filtra que sean bonos de tipo consulta, utilizados y que la fecha de la consulta coincida con la fecha pedida</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target3.InsertXML($xml3)

# --- Edit 4: lone "Notas:" paragraph after REGISTRAR LLEGADA heading -------
$full = $d.Content
$full.Find.Execute("REGISTRAR LLEGADA", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rest = $d.Range($full.End, $d.Content.End)
$rest.Find.Execute("Notas:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target4 = $d.Range($rest.Start, $rest.End)

$xml4 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:lastRenderedPageBreak/><w:t>Notas:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target4.InsertXML($xml4)
